# Update the "Michaelmas 2023" subtitle on the title slide to "Hilary 2024".
#
# The title slide (sldId 256 / first slide) has a shape named "CustomShape 2"
# (shape id 153) whose second paragraph holds the term/year line as two runs:
#   Run 1: "Michaelmas"
#   Run 2: " 2023"
# These become "Hilary" and " 2024" respectively, leaving every other run
# and all other slide content untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)

$tr = $sh.TextFrame.TextRange
$termPara = $tr.Paragraphs(2)

$termPara.Runs(1).Text = "Hilary"
$termPara.Runs(2).Text = " 2024"
